$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Append " i igrače koji su online." to the paragraph that ends
#    with "...skoriju istoriju chat-a", tagging the new run with the
#    sr-Latn-RS language, then re-create the _GoBack bookmark right
#    after that new text (it previously lived in the empty paragraph
#    at the very end of the document, after the last table).
# ------------------------------------------------------------------

# Remove the pre-existing _GoBack bookmark (it will be re-added after
# the newly inserted text further below).
$goBack = $d.Bookmarks("_GoBack")
$null = $goBack.Delete()

$target = $d.Content
$null = $target.Find.Execute("skoriju istoriju chat-a", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$target.InsertAfter(" i igrače koji su online.")
$target.LanguageID = "sr-Latn-RS"
$target.Collapse(0)

# $target now sits exactly at the end of the paragraph's content
# (right before the paragraph mark). Placing a zero-length bookmark
# there directly is unreliable, so nudge past a throwaway character
# and collapse back in front of it before adding the bookmark, then
# remove the throwaway character again.
$target.InsertAfter("#")
$target.Collapse(1)
$null = $d.Bookmarks.Add("_GoBack", $target)
$placeholderStart = $target.End
$placeholder = $d.Range($placeholderStart, $placeholderStart + 1)
$null = $placeholder.Delete()

# ------------------------------------------------------------------
# 2. The paragraph that used to hold the _GoBack bookmark (right
#    after the final table, right before the sectPr) is now empty.
# ------------------------------------------------------------------
